$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.204747333333333
$ws.Range("H2").Value = 9.614241999999999
$ws.Range("I2").Value = 0.01973032100547387
$ws.Range("J2").Value = 0.01973032100547387
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 5.943433603225555
$ws.Range("R2").Value = 53.49090242902999
$ws.Range("S2").Value = 0.0003258209267796397
$ws.Range("T2").Value = 0.0003258209267796397
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.204747333333333
$ws.Range("H3").Value = 9.614241999999999
$ws.Range("I3").Value = 0.01973032100547387
$ws.Range("J3").Value = 0.01973032100547387
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 297.7650102873184
$ws.Range("R3").Value = 2679.885092585866
$ws.Range("S3").Value = 0.01632357288583327
$ws.Range("T3").Value = 0.01632357288583328
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.204747333333333
$ws.Range("H4").Value = 9.614241999999999
$ws.Range("I4").Value = 0.01973032100547387
$ws.Range("J4").Value = 0.01973032100547387
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 0.5196049136373334
$ws.Range("R4").Value = 4.676444222736
$ws.Range("S4").Value = 0.00002848490718036977
$ws.Range("T4").Value = 0.00002848490718036978
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.204747333333333
$ws.Range("H5").Value = 9.614241999999999
$ws.Range("I5").Value = 0.01973032100547387
$ws.Range("J5").Value = 0.01973032100547387
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 55.68085583677222
$ws.Range("R5").Value = 501.1277025309499
$ws.Range("S5").Value = 0.003052442285680589
$ws.Range("T5").Value = 0.003052442285680589
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 145.2141163333334
$ws.Range("H6").Value = 435.6423490000001
$ws.Range("I6").Value = 0.8940240311559332
$ws.Range("J6").Value = 0.8940240311559333
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 269.3099857518372
$ws.Range("R6").Value = 2423.789871766536
$ws.Range("S6").Value = 0.01476365936031559
$ws.Range("T6").Value = 0.01476365936031559
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 145.2141163333334
$ws.Range("H7").Value = 435.6423490000001
$ws.Range("I7").Value = 0.8940240311559332
$ws.Range("J7").Value = 0.8940240311559333
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 13492.38437430393
$ws.Range("R7").Value = 121431.4593687354
$ws.Range("S7").Value = 0.7396568170488237
$ws.Range("T7").Value = 0.7396568170488239
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.2141163333334
$ws.Range("H8").Value = 435.6423490000001
$ws.Range("I8").Value = 0.8940240311559332
$ws.Range("J8").Value = 0.8940240311559333
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 23.54443596582134
$ws.Range("R8").Value = 211.8999236923921
$ws.Range("S8").Value = 0.001290713493076548
$ws.Range("T8").Value = 0.001290713493076548
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.2141163333334
$ws.Range("H9").Value = 435.6423490000001
$ws.Range("I9").Value = 0.8940240311559332
$ws.Range("J9").Value = 0.8940240311559333
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 2523.021454115864
$ws.Range("R9").Value = 22707.19308704278
$ws.Range("S9").Value = 0.1383128412537173
$ws.Range("T9").Value = 0.1383128412537173
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8052786666666667
$ws.Range("H10").Value = 2.415836
$ws.Range("I10").Value = 0.004957771998726471
$ws.Range("J10").Value = 0.004957771998726472
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 1.493446998971111
$ws.Range("R10").Value = 13.44102299074
$ws.Range("S10").Value = 0.00008187124106795081
$ws.Range("T10").Value = 0.00008187124106795086
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8052786666666667
$ws.Range("H11").Value = 2.415836
$ws.Range("I11").Value = 0.004957771998726471
$ws.Range("J11").Value = 0.004957771998726472
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 74.8214400461809
$ws.Range("R11").Value = 673.3929604156281
$ws.Range("S11").Value = 0.004101735220126549
$ws.Range("T11").Value = 0.00410173522012655
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.8052786666666667
$ws.Range("H12").Value = 2.415836
$ws.Range("I12").Value = 0.004957771998726471
$ws.Range("J12").Value = 0.004957771998726472
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 0.1305646618986667
$ws.Range("R12").Value = 1.175081957088
$ws.Range("S12").Value = 0.000007157596430690614
$ws.Range("T12").Value = 0.000007157596430690616
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.8052786666666667
$ws.Range("H13").Value = 2.415836
$ws.Range("I13").Value = 0.004957771998726471
$ws.Range("J13").Value = 0.004957771998726472
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 13.99130748334445
$ws.Range("R13").Value = 125.9217673501
$ws.Range("S13").Value = 0.00076700794110128
$ws.Range("T13").Value = 0.0007670079411012801
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.203389
$ws.Range("H14").Value = 39.610167
$ws.Range("I14").Value = 0.08128787583986632
$ws.Range("J14").Value = 0.08128787583986634
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 24.48663114337834
$ws.Range("R14").Value = 220.379680290405
$ws.Range("S14").Value = 0.001342364933380739
$ws.Range("T14").Value = 0.001342364933380739
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.203389
$ws.Range("H15").Value = 39.610167
$ws.Range("I15").Value = 0.08128787583986632
$ws.Range("J15").Value = 0.08128787583986634
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 1226.776045811766
$ws.Range("R15").Value = 11040.98441230589
$ws.Range("S15").Value = 0.06725225431651584
$ws.Range("T15").Value = 0.06725225431651587
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.203389
$ws.Range("H16").Value = 39.610167
$ws.Range("I16").Value = 0.08128787583986632
$ws.Range("J16").Value = 0.08128787583986634
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 2.140744678904001
$ws.Range("R16").Value = 19.26670211013601
$ws.Range("S16").Value = 0.0001173563064455779
$ws.Range("T16").Value = 0.0001173563064455779
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.203389
$ws.Range("H17").Value = 39.610167
$ws.Range("I17").Value = 0.08128787583986632
$ws.Range("J17").Value = 0.08128787583986634
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 229.4021721522584
$ws.Range("R17").Value = 2064.619549370325
$ws.Range("S17").Value = 0.01257590028352416
$ws.Range("T17").Value = 0.01257590028352416
